$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# New row 229 ("End 8-1") - row did not exist before; copy formatting
# from row 230 (A:E) so the new cells pick up style s="16" etc.
# ------------------------------------------------------------------
$ws.Range("A230:E230").Copy()
$ws.Range("A229:E229").PasteSpecial(-4122)
$ws.Range("A229").Value = "End 8-1"
$ws.Range("B229").Value = 28372
# Row 229 only has A and B populated - clear the rest that the format
# copy may have left around so no stray empty cells remain beyond B.
$ws.Range("C229:E229").Clear()

# ------------------------------------------------------------------
# Row 230 ("Enter 8-2") - values shift by -2, formula untouched so it
# keeps its shared grouping and recalculates automatically.
# ------------------------------------------------------------------
$ws.Range("B230").Value = 28731

# ------------------------------------------------------------------
# Row 231 ("Enter pipe") - values shift by -2, E231 removed entirely.
# ------------------------------------------------------------------
$ws.Range("B231").Value = 29103
$ws.Range("E231").Clear()

# ------------------------------------------------------------------
# Row 232 ("Enter pipe") - values shift by -2, E232 removed entirely.
# ------------------------------------------------------------------
$ws.Range("B232").Value = 29836
$ws.Range("E232").Clear()

# ------------------------------------------------------------------
# Row 233 ("Enter small pipe") - values shift by -2.
# ------------------------------------------------------------------
$ws.Range("B233").Value = 30492

# ------------------------------------------------------------------
# Row 234 ("Enter pipe") - values shift by -2.
# ------------------------------------------------------------------
$ws.Range("B234").Value = 30764

# ------------------------------------------------------------------
# Row 235 ("Enter pipe") - values shift by -2.
# ------------------------------------------------------------------
$ws.Range("B235").Value = 31041

# ------------------------------------------------------------------
# Row 236 ("Enter pipe (wall jump)") - values shift by -2.
# ------------------------------------------------------------------
$ws.Range("B236").Value = 31401

# ------------------------------------------------------------------
# Row 237 ("Enter Top pipe") - new C237, D237 switches from the old
# "E237-B237" one-off formula to the shared "C237-B237" pattern, and
# E237 is removed entirely.
# ------------------------------------------------------------------
$ws.Range("B237").Value = 31876
$ws.Range("C237").Value = 32382
$ws.Range("D237").Formula = "=IF(B237 >  0,C237-B237, 0)"
$ws.Range("E237").Clear()

# ------------------------------------------------------------------
# Row 238 ("End Level") - new C238 (D238 already referenced C238-B238
# via its shared formula, so it recalculates on its own). E238 removed.
# ------------------------------------------------------------------
$ws.Range("B238").Value = 32658
$ws.Range("C238").Value = 33164
$ws.Range("E238").Clear()

# ------------------------------------------------------------------
# Row 239 - previously just a stray shared-formula cell in D239; now a
# full row ("Enter 8-T").
# ------------------------------------------------------------------
$ws.Range("A230:C230").Copy()
$ws.Range("A239:C239").PasteSpecial(-4122)
$ws.Range("A239").Value = "Enter 8-T"
$ws.Range("B239").Value = 33093
$ws.Range("C239").Value = 33599
$ws.Range("D239").Formula = "=IF(B239 >  0,C239-B239, 0)"

# ------------------------------------------------------------------
# Row 240 - brand new row ("Enter Bowser Jr. Room") with a text note in
# E240 and a bare value in G240 (no explicit style, like the sheet's
# other G-column notes).
# ------------------------------------------------------------------
$ws.Range("A230:E230").Copy()
$ws.Range("A240:E240").PasteSpecial(-4122)
$ws.Range("A240").Value = "Enter Bowser Jr. Room"
$ws.Range("B240").Value = 35071
$ws.Range("C240").Value = 35593
$ws.Range("D240").Formula = "=IF(B240 >  0,C240-B240, 0)"
$ws.Range("E240").Value = "37150ish"
$ws.Range("G240").Value = 36591

# ------------------------------------------------------------------
# Row 241 - previously held only a leftover B241/E241/F241 fragment;
# now a full "End Level" row. Drop E241, add A/C/D.
# ------------------------------------------------------------------
$ws.Range("A230:D230").Copy()
$ws.Range("A241:D241").PasteSpecial(-4122)
$ws.Range("A241").Value = "End Level"
$ws.Range("B241").Value = 35976
$ws.Range("C241").Value = 36498
$ws.Range("D241").Formula = "=IF(B241 >  0,C241-B241, 0)"
$ws.Range("E241").Clear()

# ------------------------------------------------------------------
# Row 242 - new row ("Enter 8-3").
# ------------------------------------------------------------------
$ws.Range("A230:D230").Copy()
$ws.Range("A242:D242").PasteSpecial(-4122)
$ws.Range("A242").Value = "Enter 8-3"
$ws.Range("B242").Value = 36666
$ws.Range("C242").Value = 37222
$ws.Range("D242").Formula = "=IF(B242 >  0,C242-B242, 0)"

# ------------------------------------------------------------------
# Row 243 - new row ("Enter Pipe").
# ------------------------------------------------------------------
$ws.Range("A230:D230").Copy()
$ws.Range("A243:D243").PasteSpecial(-4122)
$ws.Range("A243").Value = "Enter Pipe"
$ws.Range("B243").Value = 37108
$ws.Range("C243").Value = 37664
$ws.Range("D243").Formula = "=IF(B243 >  0,C243-B243, 0)"

# ------------------------------------------------------------------
# Row 244 - new row ("End Level").
# ------------------------------------------------------------------
$ws.Range("A230:D230").Copy()
$ws.Range("A244:D244").PasteSpecial(-4122)
$ws.Range("A244").Value = "End Level"
$ws.Range("B244").Value = 41966
$ws.Range("C244").Value = 42521
$ws.Range("D244").Formula = "=IF(B244 >  0,C244-B244, 0)"

# ------------------------------------------------------------------
# Row 245 - new row ("Enter 8-4").
# ------------------------------------------------------------------
$ws.Range("A230:D230").Copy()
$ws.Range("A245:D245").PasteSpecial(-4122)
$ws.Range("A245").Value = "Enter 8-4"
$ws.Range("B245").Value = 42394
$ws.Range("C245").Value = 42949
$ws.Range("D245").Formula = "=IF(B245 >  0,C245-B245, 0)"

# ------------------------------------------------------------------
# Final selection/scroll position, matching the published view state.
# ------------------------------------------------------------------
$ws.Range("B246").Select()
